$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at 65 (pushes old row 65.. down by one) ---
$ws.Rows.Item(65).Insert()

$ws.Range("A65").Value = 1
$ws.Range("B65").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value = "Arica y Parinacota"
$ws.Range("D65").Value = 44754
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = 100112040
$ws.Range("G65").Value = "Cilantro"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 300
$ws.Range("K65").Value = 3500
$ws.Range("L65").Value = 4000
$ws.Range("M65").Value = 3750
$ws.Range("N65").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 1875
$ws.Range("Q65").Value = 2
$ws.Range("R65").Value = "Hortaliza"

# --- Insert another new row at 67 (old row 65's shifted copy is now at 66; push 66.. down) ---
$ws.Rows.Item(67).Insert()

$ws.Range("A67").Value = 1
$ws.Range("B67").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C67").Value = "Arica y Parinacota"
$ws.Range("D67").Value = 44260
$ws.Range("E67").Value = 15
$ws.Range("F67").Value = 100112040
$ws.Range("G67").Value = "Cilantro"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 200
$ws.Range("K67").Value = 4800
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = 4900
$ws.Range("N67").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O67").Value = "Región de Arica y Parinacota"
$ws.Range("P67").Value = 2450
$ws.Range("Q67").Value = 2
$ws.Range("R67").Value = "Hortaliza"
